# Generate Report for Handback
# Adds a new "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" / "Error Detail" entry for the f2863696-630e-4398-b25c-cf2909eb62eb
# row (row 8) on both the zh-cn and de-de sheets, widens column P (Error
# Detail) so the message is readable, and links the new "Latest Target File"
# cell back to the handback markdown file (same target as column A's link).

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c85b652eecd753b602e886ec6b2ccd5c07d0a00b/e2e/f2863696-630e-4398-b25c-cf2909eb62eb.md"
$mdName = "f2863696-630e-4398-b25c-cf2909eb62eb.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/040513861440e01d3e09f7f2c77a9f70e8f3400a/e2e/f2863696-630e-4398-b25c-cf2909eb62eb.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c85b652eecd753b602e886ec6b2ccd5c07d0a00b/e2e/f2863696-630e-4398-b25c-cf2909eb62eb.md."

# Column P (Error Detail) needs to be a lot wider to show the new message.
$newColumnWidth = 39.166666666666664

function Update-LocalizationSheet {
    param(
        [string]$SheetName,
        [string]$XliffName,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Widen the Error Detail column (P) to fit the new message.
    $ws.Columns.Item(16).ColumnWidth = $newColumnWidth

    # I8: Latest Target File -> hyperlinked handback markdown file name.
    $ws.Hyperlinks.Add($ws.Range("I8"), $targetUrl, "", "", $mdName)

    # J8: Latest Handback File -> the xliff produced for the handback.
    $ws.Range("J8").Value = $XliffName

    # K8: Latest Handback DateTime -> when the (out of date) handback came in.
    $ws.Range("K8").Value = $HandbackDateTime

    # P8: Error Detail -> explains the handback file is stale.
    $ws.Range("P8").Value = $errorDetail
}

Update-LocalizationSheet "zh-cn" "f2863696-630e-4398-b25c-cf2909eb62eb.458ab828df4b79f4e6ee3f905c915c303d2b9106.zh-cn.xlf" "2016-09-02 06:50:47"
Update-LocalizationSheet "de-de" "f2863696-630e-4398-b25c-cf2909eb62eb.458ab828df4b79f4e6ee3f905c915c303d2b9106.de-de.xlf" "2016-09-02 06:50:54"
